$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 73 (shifts old 73->75, 74->76, 75->77)
$ws.Rows("73:74").Insert()

# Re-create the merged cells for the two new rows (Insert does not replicate merges)
$ws.Range("A73:C73").Merge()
$ws.Range("D73:E73").Merge()
$ws.Range("F73:H73").Merge()
$ws.Range("L73:M73").Merge()
$ws.Range("A74:C74").Merge()
$ws.Range("D74:E74").Merge()
$ws.Range("F74:H74").Merge()
$ws.Range("L74:M74").Merge()

# Row 73 - item 8257
$ws.Cells.Item(73, 1).Value = 8257
$ws.Cells.Item(73, 4).Value = "AUTOMATION RV"
$ws.Cells.Item(73, 6).Value = " Part 1"
$ws.Cells.Item(73, 9).Value = " MN01"
$ws.Cells.Item(73, 10).Value = " 04/01/2032"
$ws.Cells.Item(73, 11).Value = " LB01"
$ws.Cells.Item(73, 12).Value = " XX"
$ws.Cells.Item(73, 14).Value = 1
$ws.Cells.Item(73, 15).Value = "  DEFAULTBIN"
$ws.Cells.Item(73, 16).Value = " 33159767"
$ws.Cells.Item(73, 17).Value = " 06/17/2022"
$ws.Cells.Item(73, 18).Value = "automation         06/20/2022"
$ws.Cells.Item(73, 21).Value = "WOI - 1688912"
$ws.Cells.Item(73, 22).Value = " Cycle Count"

# Row 74 - item 8258
$ws.Cells.Item(74, 1).Value = 8258
$ws.Cells.Item(74, 4).Value = "AUTOMATION RV"
$ws.Cells.Item(74, 6).Value = " Part 3"
$ws.Cells.Item(74, 9).Value = " MN03"
$ws.Cells.Item(74, 10).Value = " 04/01/2031"
$ws.Cells.Item(74, 11).Value = " LB03"
$ws.Cells.Item(74, 12).Value = " SERIAL1655440437263"
$ws.Cells.Item(74, 14).Value = 1
$ws.Cells.Item(74, 15).Value = "  AUTOMATION"
$ws.Cells.Item(74, 16).Value = " 33159750"
$ws.Cells.Item(74, 17).Value = " 06/16/2022"
$ws.Cells.Item(74, 18).Value = "automation         06/20/2022"
$ws.Cells.Item(74, 21).Value = "WOI - 1688912"
$ws.Cells.Item(74, 22).Value = " Cycle Count"

# Update the SubTotal / Total rows (now shifted to 76 / 77) from 61 to 63
$ws.Cells.Item(76, 14).Value = 63
$ws.Cells.Item(77, 14).Value = 63

# Update footer print date
$ws.PageSetup.OddFooter = "&L&`"Arial,Regular`"&9 *Data provided should be considered CONFIDENTIAL and PROPRIETARY* `n&`"-,Regular`"(c)2011-2022 MNX Global Logistics. All rights reserved. `n&`"-,Regular`"6/19/2022 9:46 PM `n&`"-,Regular`"Print User : automation &R&`"Arial,Regular`"&9Page &P of &N"
